$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "ine ilegible"
$ws.Range("E2").Value = "validado"
$ws.Range("I2").Value = "12/02/2020 09:31 a. m.;12/02/2020 09:31 a. m.;"
$ws.Range("J2").Value = "00185537"
